# "variacion entre 3 y 10"
# Update the "grilla de pruebas" sheet inputs: switch the position to BUY,
# update the capital (B3) and entry price (F3), and fill down the
# formulas in row 9 (E9/F9) following the existing pattern from rows 3-8.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grilla de pruebas")
$ws.Activate()

$ws.Range("B1").Value = "BUY"
$ws.Range("B3").Value = 1754.18
$ws.Range("F3").Value = 3.144

$ws.Range("E9").Formula = "=E8*(1+`$B`$6/100)"
$ws.Range("F9").Formula = "=IF(`$B`$1=""BUY"",F8*(1-`$B`$5/100),F8*(1+`$B`$5/100))"

[void]$ws.Range("B1").Select()
